$d = $word.ActiveDocument

function FindParaIndex($pattern) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text -match $pattern) {
            return $i
        }
    }
    return -1
}

function ReplaceParaText($para, $newText) {
    $r = $para.Range
    $old = $r.Text.TrimEnd([char]13)
    $ok = $r.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    return $ok
}

$ndash = [char]0x2013
$rsquo = [char]0x2019

# ------------------------------------------------------------------
# Locate the "Add product view" paragraph (the one that gets reworked)
# ------------------------------------------------------------------
$addIdx = FindParaIndex("^Add product view")

# The paragraph immediately before it is one of three consecutive empty
# paragraphs; delete it (it collapses away entirely).
$d.Paragraphs.Item($addIdx - 1).Range.Delete()
$addIdx = $addIdx - 1

# The (now) previous paragraph becomes the new "17-JAN-23" Heading 3.
$headingPara = $d.Paragraphs.Item($addIdx - 1)
$headingPara.Range.Text = "17-JAN-23"
$headingPara.Range.set_Style("Heading 3")

# ------------------------------------------------------------------
# Rewrite the "Add product view ..." paragraph text
# ------------------------------------------------------------------
$addPara = $d.Paragraphs.Item($addIdx)
$newAddText = "Added product view as razor page $ndash this is a bit clunky going back and forward between app and page. The only advantage is that the product page is accessible by crawlers."
ReplaceParaText $addPara $newAddText | Out-Null

# ------------------------------------------------------------------
# Fill in the bookmark-only paragraph that follows
# ------------------------------------------------------------------
$bmParaIdx = $addIdx + 1
$bmPara = $d.Paragraphs.Item($bmParaIdx)
$pStart = $bmPara.Range.Start
$pEnd = $bmPara.Range.End

# Text that must appear *after* the bookmark (a single trailing space run).
$afterPos = $pEnd - 1
$afterRange = $d.Range($afterPos, $afterPos)
$afterRange.InsertAfter(" ")

# Text that must appear *before* the bookmark.
$newBmText = "Added razor page category page. This one will be easier to crawl as it will go through all the content. Categories -> Category -> Products -> Product. The pages are not complex $ndash it is just a demo site. It" + $rsquo + "s the functionality which is important. Server-side rendering is still an important part of web for seo and non-javascript enabled browsers."
$beforeRange = $d.Range($pStart, $pStart)
$beforeRange.InsertBefore($newBmText)

# ------------------------------------------------------------------
# Delete the trailing empty paragraph that sits between the bookmark
# paragraph and "Steps:"
# ------------------------------------------------------------------
$stepsIdx = FindParaIndex("^Steps:")
$d.Paragraphs.Item($stepsIdx - 1).Range.Delete()
